$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Font.TotallyBogusPropertyXYZ = 123
Write-Output "done"
